# Add a new "2022" column (S) to the right of the existing "2021" column (R),
# carrying over R's per-row formatting (font, borders, number format, alignment)
# onto the new S cells, then move the active selection to T6 (matching the
# author's final cursor position after the new column was added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newColumnValues = [ordered]@{
    4  = 2022
    5  = 4.9538761752705343
    6  = 11.304954640614097
    7  = 5.1593323216995444
    8  = 13.687943262411348
    9  = 10.22864019253911
    10 = 9.1213700670141478
    11 = 3.1335149863760217
    12 = 2.872905173311127
    13 = 3.527842284697861
    14 = 5.0305321314335565
}

foreach ($row in $newColumnValues.Keys) {
    $source = $ws.Range("R$row")
    $target = $ws.Range("S$row")

    # Copy R's formatting (font/border/numberformat/alignment) onto S...
    $source.Copy($target)
    # ...then overwrite with the new value for this column/row.
    $target.Value = $newColumnValues[$row]
}

$ws.Range("T6").Select()
